$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AM1").Value = 0.91207723270377083
$ws.Range("BM1").Value = 0.84585695280622797
$ws.Range("BP1").Value = 0.7641938340251575
$ws.Range("A2").Value = 0.61421990307814167
$ws.Range("C2").Value = 0.86420906889003479
$ws.Range("B4").Value = 0.90830103943611773
$ws.Range("C4").Value = 0.98146241393940137
$ws.Range("E4").Value = 0.93465410260925785
$ws.Range("BO4").Value = 0.73793214392535988
$ws.Range("C5").Value = 0.80136870706702823
$ws.Range("F5").Value = 0.87156438617613929
$ws.Range("N6").Value = 0.9972380637559124
$ws.Range("E7").Value = 0.61895178425831632
$ws.Range("F7").Value = 0.89960184831337064
$ws.Range("I7").Value = 0.87281389961957312
$ws.Range("G8").Value = 0.82707861019934814
$ws.Range("I8").Value = 0.76094905744589547
$ws.Range("BA8").Value = 0.92131325770011241
$ws.Range("J9").Value = 0.72131183503423735
$ws.Range("K9").Value = 0.88232729915159047
$ws.Range("O10").Value = 0.95293717377476161
$ws.Range("J11").Value = 0.98608819431701589
$ws.Range("M11").Value = 0.92349946035028219
$ws.Range("J12").Value = 0.73128215729039181
$ws.Range("K12").Value = 0.83994372470714618
$ws.Range("BJ12").Value = 0.98477671680416479
$ws.Range("L13").Value = 0.69230276786967182
$ws.Range("N13").Value = 0.95069560657070107
$ws.Range("BA13").Value = 0.83929361797186297
$ws.Range("P14").Value = 0.79508691945683507
$ws.Range("M15").Value = 0.90468213643188566
$ws.Range("N15").Value = 0.70502926210401873
$ws.Range("Z15").Value = 0.82751209928013481
$ws.Range("R16").Value = 0.89172492603551134
$ws.Range("P17").Value = 0.97153779891775516
$ws.Range("Q18").Value = 0.89277813937436701
$ws.Range("Q19").Value = 0.82595360544068663
$ws.Range("R19").Value = 0.72216897486564513
$ws.Range("T19").Value = 0.88699716393286177
$ws.Range("U19").Value = 0.72622316098563311
$ws.Range("L20").Value = 0.93031775678043394
$ws.Range("R20").Value = 0.82606071714919538
$ws.Range("V20").Value = 0.70805776168685719
$ws.Range("D21").Value = 0.60661264721659136
$ws.Range("T21").Value = 0.84190048220305858
$ws.Range("U23").Value = 0.88135783774130549
$ws.Range("V23").Value = 0.86088596908055803
$ws.Range("Y23").Value = 0.90130560059446463
$ws.Range("V24").Value = 0.73133663313606867
$ws.Range("W24").Value = 0.80820376435728847
$ws.Range("Y24").Value = 0.99739590685182111
$ws.Range("Z25").Value = 0.70761930348990942
$ws.Range("X26").Value = 0.92237537760980737
$ws.Range("BP26").Value = 0.63514680505638643
$ws.Range("Y27").Value = 0.75330945699103558
$ws.Range("Z27").Value = 0.97997326339831414
$ws.Range("AC27").Value = 0.79630293788529483
$ws.Range("R28").Value = 0.99031309130410627
$ws.Range("AA28").Value = 0.89345296218216619
$ws.Range("AD28").Value = 0.68223902156791782
$ws.Range("AD29").Value = 0.97354279412868627
$ws.Range("AS29").Value = 0.82453856559266203
$ws.Range("AE30").Value = 0.7822177638907426
$ws.Range("AF30").Value = 0.67139363662681206
$ws.Range("AF31").Value = 0.65553528334124345
$ws.Range("AG31").Value = 0.92803643494449151
$ws.Range("AG32").Value = 0.76530985462116985
$ws.Range("AH33").Value = 0.92439249344571572
$ws.Range("AI33").Value = 0.813721992789602
$ws.Range("AF34").Value = 0.79662372024757044
$ws.Range("Y35").Value = 0.86149170568504285
$ws.Range("AH35").Value = 0.85358708561320407
$ws.Range("BD35").Value = 0.9650919400292941
$ws.Range("AH36").Value = 0.97670505172092792
$ws.Range("BH36").Value = 0.94572242805476259
$ws.Range("AJ37").Value = 0.74324140698596119
$ws.Range("AL37").Value = 0.85459500218770379
$ws.Range("AM37").Value = 0.84268981671079968
$ws.Range("BE37").Value = 0.85542970316889733
$ws.Range("AN38").Value = 0.69498478526588259
$ws.Range("BP38").Value = 0.9787522658676282
$ws.Range("AG39").Value = 0.89456117994205209
$ws.Range("AO40").Value = 0.96554568291360066
$ws.Range("AP40").Value = 0.71471198051966534
$ws.Range("AO42").Value = 0.98041259761070343
$ws.Range("AQ42").Value = 0.77698063169542997
$ws.Range("AU42").Value = 0.76478472193678049
$ws.Range("BB42").Value = 0.95076580913580788
$ws.Range("BJ42").Value = 0.95522869597772808
$ws.Range("AO43").Value = 0.823127290856706
$ws.Range("AS43").Value = 0.85574572542465077
$ws.Range("AP44").Value = 0.82807666247603473
$ws.Range("AQ44").Value = 0.78299406345503086
$ws.Range("AS44").Value = 0.96362102581137909
$ws.Range("AT44").Value = 0.92699137092526762
$ws.Range("AT45").Value = 0.94248361230242994
$ws.Range("AV46").Value = 0.63411160145046086
$ws.Range("AS47").Value = 0.76937192501929275
$ws.Range("AT47").Value = 0.98677564083008718
$ws.Range("AV47").Value = 0.68279524820257287
$ws.Range("AK48").Value = 0.8539461082645311
$ws.Range("AV49").Value = 0.86060584509746307
$ws.Range("AY49").Value = 0.71341556693495756
$ws.Range("AW50").Value = 0.6429951188603158
$ws.Range("AY50").Value = 0.82694533720475882
$ws.Range("AZ50").Value = 0.70802970914185803
$ws.Range("A51").Value = 0.64123794914225773
$ws.Range("AY52").Value = 0.85157023808742882
$ws.Range("BA52").Value = 0.70047811700970386
$ws.Range("BK52").Value = 0.51720260032968501
$ws.Range("AY53").Value = 0.86942522457044213
$ws.Range("BA54").Value = 0.74567687142036942
$ws.Range("BD54").Value = 0.93673785290522349
$ws.Range("BD55").Value = 0.6312430020539006
$ws.Range("BC57").Value = 0.86959165422353268
$ws.Range("BD57").Value = 0.87521823345923877
$ws.Range("BF57").Value = 0.89877788693486793
$ws.Range("BD58").Value = 0.88268205564619628
$ws.Range("BG58").Value = 0.94489083890559711
$ws.Range("V59").Value = 0.86819037827862955
$ws.Range("AQ60").Value = 0.8743955675273688
$ws.Range("BF60").Value = 0.97027991569482386
$ws.Range("BG60").Value = 0.70890245452385325
$ws.Range("BI60").Value = 0.9930800041197998
$ws.Range("BJ61").Value = 0.60239491834565095
$ws.Range("BK61").Value = 0.89826722009833193
$ws.Range("BK64").Value = 0.96785843958027629
$ws.Range("BK65").Value = 0.88282375501624277
$ws.Range("BL65").Value = 0.93260502049589133
$ws.Range("BO65").Value = 0.80653593475976926
$ws.Range("BL66").Value = 0.8141706126295376
$ws.Range("BO66").Value = 0.55068170762639546
$ws.Range("BP66").Value = 0.88820536520189997
$ws.Range("F67").Value = 0.78330067946174653
$ws.Range("BP67").Value = 0.90898041826470211
$ws.Range("C68").Value = 0.76133676504319614
